$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Turn off automatic calculation (calcMode="manual")
$excel.Calculation = -4135   # xlCalculationManual

# Add the new data row (row 33)
$ws.Range("A33").Value = 10002
$ws.Range("B33").Value = 110032
$ws.Range("C33").Value = 10032
$ws.Range("D33").Value = "eng"
$ws.Range("E33").Value = $true
$ws.Range("F33").Value = "superadmin"
$ws.Range("G33").Value = "now()"
$ws.Range("H33").Value = "now()"

# Update selection / view state to match end result
$ws.Range("E31").Select()
